$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F entirely (header + data) - it no longer exists in the target layout
$ws.Columns.Item(6).Delete()

# Rewrite the data rows (2..28) with the updated dataset
$data = @(
    @("2021-10", -1.9, 753277, 7409924, 61.7),
    @("2021-11", 15.8, 858566, 8346016, 57.3),
    @("2021-12", 2.2, 902551, 9214377, 48.9),
    @("2021-02", "", "", 956287, 132.3),
    @("2021-03", 91.1, 748288, 1695764, 109.1),
    @("2021-04", 100.1, 703019, 2391928, 106),
    @("2021-05", 49.2, 700247, 3092193, 89.9),
    @("2021-06", 67.9, 817251, 3925950, 85.5),
    @("2021-07", 18.4, 596110, 4834380, 84.6),
    @("2021-08", 38.7, 680220, 5894317, 89.4),
    @("2021-09", 4.1, 701306, 6595643, 73.6),
    @("2022-10", -30.1, 534035, 5427892, -27),
    @("2022-11", -38, 538496, 5967043, -28.8),
    @("2022-12", -46.5, 490992, 6458014, -30.3),
    @("2022-02", "", "", 1079696, -2),
    @("2022-03", -23.7, 694150, 1785543, -10.9),
    @("2022-04", -43, 482336, 2265142, -20.3),
    @("2022-05", -33.2, 529558, 2793997, -23.2),
    @("2022-06", -48.1, 476557, 3277196, -28.3),
    @("2022-07", -33.5, 467654, 3748735, -28.9),
    @("2022-08", -37.4, 426813, 4176153, -29.5),
    @("2022-09", -3.8, 679899, 4877858, -26.4),
    @("2023-02", "", "", 755968, -28.7),
    @("2023-03", -1.7, 700272, 1450428, -20.6),
    @("2023-04", 47.6, 706001, 2178790, -2.3),
    @("2023-05", 34.3, 697262, 2867635, 4.5),
    @("2023-06", 40, 662658, 3530430, 9.6),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    if ($vals[1] -eq "") {
        $ws.Cells.Item($row, 2).Value = ""
    } else {
        $ws.Cells.Item($row, 2).Value = $vals[1]
    }
    if ($vals[2] -eq "") {
        $ws.Cells.Item($row, 3).Value = ""
    } else {
        $ws.Cells.Item($row, 3).Value = $vals[2]
    }
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Column A (the period labels) carries the bold/bordered/centered header style;
# make sure every data row down to the new last row (28) has it too.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
